$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$c = $t.Cell(1, 1)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "827÷3=275, 2") {
    Write-Host "MISMATCH at row 1 col 1: " $c.Range.Text
}
$c.Range.Text = "697÷5=139, 2"

$c = $t.Cell(1, 2)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "402÷3=134, 0") {
    Write-Host "MISMATCH at row 1 col 2: " $c.Range.Text
}
$c.Range.Text = "520÷2=260, 0"

$c = $t.Cell(1, 3)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "152÷6=25, 2") {
    Write-Host "MISMATCH at row 1 col 3: " $c.Range.Text
}
$c.Range.Text = "336÷8=42, 0"

$c = $t.Cell(1, 4)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "376÷8=47, 0") {
    Write-Host "MISMATCH at row 1 col 4: " $c.Range.Text
}
$c.Range.Text = "778÷7=111, 1"

$c = $t.Cell(1, 5)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "974÷5=194, 4") {
    Write-Host "MISMATCH at row 1 col 5: " $c.Range.Text
}
$c.Range.Text = "318÷6=53, 0"

$c = $t.Cell(5, 1)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "515÷4=128, 3") {
    Write-Host "MISMATCH at row 5 col 1: " $c.Range.Text
}
$c.Range.Text = "854÷9=94, 8"

$c = $t.Cell(5, 2)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "262÷8=32, 6") {
    Write-Host "MISMATCH at row 5 col 2: " $c.Range.Text
}
$c.Range.Text = "505÷6=84, 1"

$c = $t.Cell(5, 3)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "878÷8=109, 6") {
    Write-Host "MISMATCH at row 5 col 3: " $c.Range.Text
}
$c.Range.Text = "133÷7=19, 0"

$c = $t.Cell(5, 4)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "891÷7=127, 2") {
    Write-Host "MISMATCH at row 5 col 4: " $c.Range.Text
}
$c.Range.Text = "119÷5=23, 4"

$c = $t.Cell(5, 5)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "790÷5=158, 0") {
    Write-Host "MISMATCH at row 5 col 5: " $c.Range.Text
}
$c.Range.Text = "699÷6=116, 3"

$c = $t.Cell(9, 1)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "132÷7=18, 6") {
    Write-Host "MISMATCH at row 9 col 1: " $c.Range.Text
}
$c.Range.Text = "777÷7=111, 0"

$c = $t.Cell(9, 2)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "243÷3=81, 0") {
    Write-Host "MISMATCH at row 9 col 2: " $c.Range.Text
}
$c.Range.Text = "140÷3=46, 2"

$c = $t.Cell(9, 3)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "433÷3=144, 1") {
    Write-Host "MISMATCH at row 9 col 3: " $c.Range.Text
}
$c.Range.Text = "595÷7=85, 0"

$c = $t.Cell(9, 4)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "198÷6=33, 0") {
    Write-Host "MISMATCH at row 9 col 4: " $c.Range.Text
}
$c.Range.Text = "224÷4=56, 0"

$c = $t.Cell(9, 5)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "221÷5=44, 1") {
    Write-Host "MISMATCH at row 9 col 5: " $c.Range.Text
}
$c.Range.Text = "265÷5=53, 0"

$c = $t.Cell(13, 1)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "921÷8=115, 1") {
    Write-Host "MISMATCH at row 13 col 1: " $c.Range.Text
}
$c.Range.Text = "623÷3=207, 2"

$c = $t.Cell(13, 2)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "844÷5=168, 4") {
    Write-Host "MISMATCH at row 13 col 2: " $c.Range.Text
}
$c.Range.Text = "829÷4=207, 1"

$c = $t.Cell(13, 3)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "516÷4=129, 0") {
    Write-Host "MISMATCH at row 13 col 3: " $c.Range.Text
}
$c.Range.Text = "145÷6=24, 1"

$c = $t.Cell(13, 4)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "184÷4=46, 0") {
    Write-Host "MISMATCH at row 13 col 4: " $c.Range.Text
}
$c.Range.Text = "857÷9=95, 2"

$c = $t.Cell(13, 5)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "915÷6=152, 3") {
    Write-Host "MISMATCH at row 13 col 5: " $c.Range.Text
}
$c.Range.Text = "831÷9=92, 3"

$c = $t.Cell(17, 1)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "991÷5=198, 1") {
    Write-Host "MISMATCH at row 17 col 1: " $c.Range.Text
}
$c.Range.Text = "586÷8=73, 2"

$c = $t.Cell(17, 2)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "191÷4=47, 3") {
    Write-Host "MISMATCH at row 17 col 2: " $c.Range.Text
}
$c.Range.Text = "772÷8=96, 4"

$c = $t.Cell(17, 3)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "316÷8=39, 4") {
    Write-Host "MISMATCH at row 17 col 3: " $c.Range.Text
}
$c.Range.Text = "374÷2=187, 0"

$c = $t.Cell(17, 4)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "878÷8=109, 6") {
    Write-Host "MISMATCH at row 17 col 4: " $c.Range.Text
}
$c.Range.Text = "790÷9=87, 7"

$c = $t.Cell(17, 5)
if ($c.Range.Text.TrimEnd([char]13, [char]7) -ne "305÷9=33, 8") {
    Write-Host "MISMATCH at row 17 col 5: " $c.Range.Text
}
$c.Range.Text = "105÷8=13, 1"
